$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.722.62'
$ws.Cells.Item(2, 5).Value = '  +0.37%  '
$ws.Cells.Item(3, 4).Value = '1.850.25'
$ws.Cells.Item(3, 5).Value = '  +0.36%  '
$ws.Cells.Item(4, 5).Value = '  +0.52%  '
$ws.Cells.Item(5, 4).Value = '''313.12'
$ws.Cells.Item(5, 5).Value = '  -0.69%  '
$ws.Cells.Item(6, 5).Value = '  +0.37%  '
$ws.Cells.Item(7, 4).Value = '''0.4282'
$ws.Cells.Item(7, 5).Value = '  +0.99%  '
$ws.Cells.Item(8, 4).Value = '''0.3589'
$ws.Cells.Item(8, 5).Value = '  -1.36%  '
$ws.Cells.Item(9, 4).Value = '''0.07314'
$ws.Cells.Item(9, 5).Value = '  +0.46%  '
$ws.Cells.Item(10, 4).Value = '''0.8715'
$ws.Cells.Item(10, 5).Value = '  -1.92%  '
$ws.Cells.Item(11, 4).Value = '''20.75'
$ws.Cells.Item(11, 5).Value = '  +0.29%  '
$ws.Cells.Item(12, 4).Value = '1.862.89'
$ws.Cells.Item(12, 5).Value = '  +3.28%  '
$ws.Cells.Item(13, 4).Value = '''6.554'
$ws.Cells.Item(13, 5).Value = '  -0.03%  '
$ws.Cells.Item(14, 4).Value = '''5.341'
$ws.Cells.Item(14, 5).Value = '  -0.14%  '
$ws.Cells.Item(15, 4).Value = '''0.06998'
$ws.Cells.Item(15, 5).Value = '  +1.75%  '
$ws.Cells.Item(16, 4).Value = '''1.007'
$ws.Cells.Item(16, 5).Value = '  +0.43%  '
$ws.Cells.Item(17, 4).Value = '''79.71'
$ws.Cells.Item(17, 5).Value = '  +0.88%  '
$ws.Cells.Item(18, 4).Value = '''0.000008965'
$ws.Cells.Item(18, 5).Value = '  +1.06%  '
$ws.Cells.Item(19, 4).Value = '''1.005'
$ws.Cells.Item(19, 5).Value = '  +0.47%  '
$ws.Cells.Item(20, 5).Value = '  -0.83%  '
$ws.Cells.Item(21, 4).Value = '27.871.39'
$ws.Cells.Item(21, 5).Value = '  +0.95%  '
$ws.Cells.Item(22, 4).Value = '''5.000'
$ws.Cells.Item(22, 5).Value = '  +0.32%  '
$ws.Cells.Item(23, 5).Value = '  -1.69%  '
$ws.Cells.Item(24, 4).Value = '2.109.61'
$ws.Cells.Item(24, 5).Value = '  +3.55%  '
$ws.Cells.Item(25, 4).Value = '''1.989'
$ws.Cells.Item(25, 5).Value = '  +2.70%  '
$ws.Cells.Item(26, 4).Value = '''155.64'
$ws.Cells.Item(26, 5).Value = '  +0.58%  '
$ws.Cells.Item(27, 4).Value = '''18.52'
$ws.Cells.Item(27, 5).Value = '  -2.14%  '
$ws.Cells.Item(28, 5).Value = '  -0.92%  '
$ws.Cells.Item(29, 4).Value = '''5.276'
$ws.Cells.Item(29, 5).Value = '  -0.02%  '
$ws.Cells.Item(30, 4).Value = '''1.872'
$ws.Cells.Item(30, 5).Value = '  +0.20%  '
$ws.Cells.Item(31, 4).Value = '''0.08925'
$ws.Cells.Item(31, 5).Value = '  -0.10%  '
$ws.Cells.Item(32, 4).Value = '''0.7655'
$ws.Cells.Item(32, 5).Value = '  -0.90%  '
$ws.Cells.Item(33, 4).Value = '''2.972'
$ws.Cells.Item(33, 5).Value = '  +1.52%  '
$ws.Cells.Item(34, 4).Value = '''4.513'
$ws.Cells.Item(34, 5).Value = '  -1.56%  '
$ws.Cells.Item(35, 4).Value = '''1.128'
$ws.Cells.Item(35, 5).Value = '  +2.84%  '
$ws.Cells.Item(36, 5).Value = '  +0.40%  '
$ws.Cells.Item(37, 4).Value = '''0.05442'
$ws.Cells.Item(37, 5).Value = '  +1.39%  '
$ws.Cells.Item(38, 4).Value = '''1.106'
$ws.Cells.Item(38, 5).Value = '  +0.59%  '
$ws.Cells.Item(39, 5).Value = '  -0.08%  '
$ws.Cells.Item(40, 4).Value = '''2.834'
$ws.Cells.Item(40, 5).Value = '  +1.16%  '
$ws.Cells.Item(41, 4).Value = '''0.1668'
$ws.Cells.Item(41, 5).Value = '  +0.95%  '
$ws.Cells.Item(42, 4).Value = '''0.5081'
$ws.Cells.Item(42, 5).Value = '  -0.18%  '
$ws.Cells.Item(43, 4).Value = '''6.614'
$ws.Cells.Item(43, 5).Value = '  -3.86%  '
$ws.Cells.Item(44, 4).Value = '''8.413'
$ws.Cells.Item(44, 5).Value = '  +1.70%  '
$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(45, 4).Value = '''106.59'
$ws.Cells.Item(45, 5).Value = '  +1.64%  '
$ws.Cells.Item(46, 2).Value = 'Cronos'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(46, 4).Value = '''0.06549'
$ws.Cells.Item(46, 5).Value = '  -0.46%  '
$ws.Cells.Item(47, 4).Value = '''10.40'
$ws.Cells.Item(47, 5).Value = '  +0.38%  '
$ws.Cells.Item(48, 4).Value = '''0.4659'
$ws.Cells.Item(48, 5).Value = '  -1.31%  '
$ws.Cells.Item(49, 5).Value = '  +0.42%  '
$ws.Cells.Item(50, 4).Value = '''1.636'
$ws.Cells.Item(50, 5).Value = '  +0.08%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Value = '''64.62'
$ws.Cells.Item(51, 5).Value = '  +0.28%  '
